# Update the "Runmode" flags for the AddCustomerTest / OpenAccountTest
# scenarios (common utility for runmodes) and leave the selection/active
# sheet the way the author left the workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestCases")
$ws2 = $wb.Worksheets.Item("TestData")

# TestCases sheet: flip OpenAccountTest's Runmode from N to Y
$ws1.Range("B3").Value = "Y"

# TestData sheet:
#  - AddCustomerTest's second data row (jyoti) Runmode: N -> Y
#  - OpenAccountTest's two data rows (manish k / jyoti k) Runmode: Y -> N
$ws2.Range("A4").Value = "Y"
$ws2.Range("A8").Value = "N"
$ws2.Range("A9").Value = "N"

# Restore the saved selection/active sheet state recorded in the workbook
$ws2.Range("A4").Select()
$ws1.Activate()
$ws1.Range("B4").Select()
